$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update position (column C) values: "软件开发" split into 后端开发/前端开发/移动端开发;
# "软件测试" split into 功能测试/性能测试
$ws.Range("C20").Value = "后端开发"
$ws.Range("C21").Value = "功能测试"
$ws.Range("C22").Value = "后端开发"
$ws.Range("C23").Value = "移动端开发"
$ws.Range("C27").Value = "后端开发"
$ws.Range("C29").Value = "后端开发"
$ws.Range("C31").Value = "后端开发"
$ws.Range("C32").Value = "后端开发"
$ws.Range("C33").Value = "后端开发"
$ws.Range("C35").Value = "功能测试"
$ws.Range("C37").Value = "功能测试"
$ws.Range("C39").Value = "功能测试"
$ws.Range("C40").Value = "后端开发"
$ws.Range("C44").Value = "后端开发"
$ws.Range("C45").Value = "后端开发"
$ws.Range("C46").Value = "后端开发"
$ws.Range("C47").Value = "功能测试"
$ws.Range("C48").Value = "性能测试"
$ws.Range("C50").Value = "功能测试"
$ws.Range("C53").Value = "后端开发"
$ws.Range("C55").Value = "后端开发"
$ws.Range("C56").Value = "后端开发"
$ws.Range("C57").Value = "后端开发"
$ws.Range("C58").Value = "后端开发"
$ws.Range("C60").Value = "后端开发"
$ws.Range("C62").Value = "前端开发"
$ws.Range("C63").Value = "后端开发"
$ws.Range("C64").Value = "功能测试"
$ws.Range("C66").Value = "前端开发"
$ws.Range("C67").Value = "功能测试"
$ws.Range("C72").Value = "后端开发"
$ws.Range("C73").Value = "移动端开发"
$ws.Range("C74").Value = "功能测试"
$ws.Range("C75").Value = "功能测试"
$ws.Range("C83").Value = "性能测试"
$ws.Range("C84").Value = "后端开发"
$ws.Range("C87").Value = "后端开发"
$ws.Range("C92").Value = "后端开发"
$ws.Range("C93").Value = "后端开发"
$ws.Range("C95").Value = "前端开发"
$ws.Range("C98").Value = "后端开发"
$ws.Range("C99").Value = "后端开发"
$ws.Range("C100").Value = "前端开发"
$ws.Range("C106").Value = "后端开发"
$ws.Range("C107").Value = "后端开发"
$ws.Range("C112").Value = "后端开发"
$ws.Range("C113").Value = "后端开发"
$ws.Range("C115").Value = "前端开发"
$ws.Range("C119").Value = "后端开发"
$ws.Range("C124").Value = "后端开发"
$ws.Range("C125").Value = "后端开发"
$ws.Range("C130").Value = "后端开发"
$ws.Range("C132").Value = "后端开发"
$ws.Range("C134").Value = "功能测试"
$ws.Range("C136").Value = "后端开发"
$ws.Range("C141").Value = "后端开发"
$ws.Range("C143").Value = "前端开发"
$ws.Range("C146").Value = "前端开发"
$ws.Range("C147").Value = "后端开发"
$ws.Range("C149").Value = "后端开发"
$ws.Range("C150").Value = "后端开发"
$ws.Range("C151").Value = "后端开发"
$ws.Range("C153").Value = "后端开发"
$ws.Range("C160").Value = "前端开发"
$ws.Range("C165").Value = "后端开发"
$ws.Range("C166").Value = "后端开发"
$ws.Range("C167").Value = "功能测试"
$ws.Range("C168").Value = "后端开发"
$ws.Range("C169").Value = "后端开发"
$ws.Range("C171").Value = "后端开发"
$ws.Range("C176").Value = "后端开发"
$ws.Range("C177").Value = "后端开发"
$ws.Range("C181").Value = "后端开发"
$ws.Range("C182").Value = "后端开发"
$ws.Range("C192").Value = "功能测试"
$ws.Range("C194").Value = "后端开发"
$ws.Range("C195").Value = "后端开发"
$ws.Range("C196").Value = "后端开发"
$ws.Range("C197").Value = "前端开发"
$ws.Range("C198").Value = "后端开发"
$ws.Range("C199").Value = "后端开发"
$ws.Range("C200").Value = "后端开发"
$ws.Range("C201").Value = "后端开发"
$ws.Range("C204").Value = "后端开发"
$ws.Range("C212").Value = "后端开发"
$ws.Range("C213").Value = "前端开发"

# Add hidden AutoFilter defined name (Excel records this when a filter range is set)
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=talents!`$A`$1:`$C`$213")
$fd.Visible = $false

# Update the active selection cell
$ws.Range("I10").Select()
